# Insert a new paragraph "Computer vision " right after the
# "Introduction" heading paragraph, matching the formatting of that
# paragraph's run (Arial, color 222222, white shading).

$d = $word.ActiveDocument

# Locate the "Introduction" paragraph.
foreach ($p in $d.Paragraphs) {
    $r = $p.Range
    $text = $r.Text.TrimEnd([char]13, [char]7)
    if ($text -eq "Introduction") {
        # Collapse to the end of this paragraph (just before its
        # paragraph mark) and insert a new paragraph after it.
        $end = $r.End - 1
        $target = $d.Range($end, $end)
        $target.InsertParagraphAfter()

        # The newly created paragraph is the one right after this one.
        $newPara = $p.Next(1)
        $newRange = $newPara.Range
        $newRange.Font.Name = "Arial"
        $newRange.Font.Color = 2236962
        $newRange.Text = "Computer vision "
        break
    }
}
